$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "challenges": shrink from 6 challenges (rows 2-7) to 4 (rows 2-5)
# ---------------------------------------------------------------------------
$wsChallenges = $wb.Worksheets.Item("challenges")

# target (L3) 12 -> 30
$wsChallenges.Cells.Item(3, 12).Value = 30

# success_next (M5) cleared - challenge 17/4 no longer has a success_next
$wsChallenges.Cells.Item(5, 13).ClearContents()

# drop the old rows 6 and 7 entirely (deleting row 6 twice removes both,
# since row 7 shifts up into row 6's place after the first delete)
$wsChallenges.Rows.Item(6).Delete()
$wsChallenges.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# Sheet "tasks": rework existing tasks + append new ones
# ---------------------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item("tasks")

# Row 2: tutorial_video(physical_activity) - only the secret token changes
$wsTasks.Cells.Item(2, 12).Value = " [SECRET, EQUAL, o5bnwon9i6rv3nupuntgzsqhywfb0sk7sppsi9mazy5bwu5ph]"

# Row 3: Do_10_push-ups -> Take_a_45-minute_walk_without_stopping
$wsTasks.Cells.Item(3, 2).Value = "Take_a_45-minute_walk_without_stopping"
$wsTasks.Cells.Item(3, 12).Value = " [SECRET, EQUAL, 60szg8o5o8]"

# Row 4: Do_10_push-ups -> Take_200_steps
$wsTasks.Cells.Item(4, 2).Value = "Take_200_steps"
$wsTasks.Cells.Item(4, 9).Value = "WALK"
$wsTasks.Cells.Item(4, 10).Value = "WALK"
$wsTasks.Cells.Item(4, 12).Value = "[STEPS, STRICTLY_GREATER, 3000], [SECRET, EQUAL, ej3dg5z2rq1m6g7v97m]"
$wsTasks.Cells.Item(4, 13).Value = 2

# Row 5: used to be tutorial_video(social_activity) (challenge 3), becomes
# Include_10_minutes_of_uphill_walking_during_one_of_your_walks (challenge 2)
$wsTasks.Cells.Item(5, 1).Value = 2
$wsTasks.Cells.Item(5, 2).Value = "Include_10_minutes_of_uphill_walking_during_one_of_your_walks"
$wsTasks.Cells.Item(5, 6).ClearContents()
$wsTasks.Cells.Item(5, 9).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(5, 10).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(5, 12).Value = " [SECRET, EQUAL, h6xgwwhtqnc2gplsy3h1ncqvt09bssdiqxi6rk8cjk9ku8u]"
$wsTasks.Cells.Item(5, 13).Value = 2

# Row 6: used to be Engage_with_others (challenge 4), becomes
# Walk_9000_steps_in_a_day (challenge 2)
$wsTasks.Cells.Item(6, 1).Value = 2
$wsTasks.Cells.Item(6, 2).Value = "Walk_9000_steps_in_a_day"
$wsTasks.Cells.Item(6, 9).Value = "DAY_AGGREGATE"
$wsTasks.Cells.Item(6, 10).Value = "DAY_AGGREGATE"
$wsTasks.Cells.Item(6, 12).Value = "[STEPS_SUM, STRICTLY_GREATER, 9000], [SECRET, EQUAL, 1aaukm7ml4g9m8]"
$wsTasks.Cells.Item(6, 13).Value = 10

# Row 7: used to be tutorial_video(cognitive_activity) (challenge 5), becomes
# another Walk_9000_steps_in_a_day (challenge 2)
$wsTasks.Cells.Item(7, 1).Value = 2
$wsTasks.Cells.Item(7, 2).Value = "Walk_9000_steps_in_a_day"
$wsTasks.Cells.Item(7, 6).ClearContents()
$wsTasks.Cells.Item(7, 9).Value = "DAY_AGGREGATE"
$wsTasks.Cells.Item(7, 10).Value = "DAY_AGGREGATE"
$wsTasks.Cells.Item(7, 12).Value = "[STEPS_SUM, STRICTLY_GREATER, 9000], [SECRET, EQUAL, qqj1vx1hq6ndxp5d6q907icqb1zeeh1t7vj81fvypxm]"
$wsTasks.Cells.Item(7, 13).Value = 10

# Row 8: used to be Play_a_boardgame (challenge 6), becomes
# tutorial_video(social_activity) (challenge 3)
$wsTasks.Cells.Item(8, 1).Value = 3
$wsTasks.Cells.Item(8, 2).Value = "tutorial_video(social_activity)"
$wsTasks.Cells.Item(8, 6).Value = "https://campaigns.healthyw8.gamebus.eu/api/media/generated-296ffd13/f0a366cc-c574-4807-8dab-5dd53dd47f70.h5p"
$wsTasks.Cells.Item(8, 9).Value = "H5P_GENERAL"
$wsTasks.Cells.Item(8, 10).Value = "H5P_GENERAL"
$wsTasks.Cells.Item(8, 12).Value = " [SECRET, EQUAL, dibhlux6stuurubsixtsfnfa5nke6v5kr]"
$wsTasks.Cells.Item(8, 13).Value = 1

# Row 9 (new): Call_a_friend (challenge 4)
$wsTasks.Cells.Item(9, 1).Value = 4
$wsTasks.Cells.Item(9, 2).Value = "Call_a_friend"
$wsTasks.Cells.Item(9, 4).Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$wsTasks.Cells.Item(9, 7).Value = 1
$wsTasks.Cells.Item(9, 8).Value = 7
$wsTasks.Cells.Item(9, 9).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(9, 10).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(9, 11).Value = 0
$wsTasks.Cells.Item(9, 12).Value = " [SECRET, EQUAL, 81hceada27ud7qcheqdudbuaqkb]"
$wsTasks.Cells.Item(9, 13).Value = 1
$wsTasks.Cells.Item(9, 14).Value = "GameBus Studio"

# Row 10 (new): Go_to_a_social_event (challenge 4)
$wsTasks.Cells.Item(10, 1).Value = 4
$wsTasks.Cells.Item(10, 2).Value = "Go_to_a_social_event"
$wsTasks.Cells.Item(10, 4).Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$wsTasks.Cells.Item(10, 7).Value = 1
$wsTasks.Cells.Item(10, 8).Value = 7
$wsTasks.Cells.Item(10, 9).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(10, 10).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(10, 11).Value = 0
$wsTasks.Cells.Item(10, 12).Value = " [SECRET, EQUAL, fa5e6ynirrcd]"
$wsTasks.Cells.Item(10, 13).Value = 2
$wsTasks.Cells.Item(10, 14).Value = "GameBus Studio"

# Row 11 (new): Call_a_friend (challenge 4)
$wsTasks.Cells.Item(11, 1).Value = 4
$wsTasks.Cells.Item(11, 2).Value = "Call_a_friend"
$wsTasks.Cells.Item(11, 4).Value = "https://campaigns.healthyw8.gamebus.eu/api/media/HW8-immutable/5ff935d3-d0ae-4dce-bfcd-d2f71bf2ca63.jpeg"
$wsTasks.Cells.Item(11, 7).Value = 1
$wsTasks.Cells.Item(11, 8).Value = 7
$wsTasks.Cells.Item(11, 9).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(11, 10).Value = "GENERAL_ACTIVITY"
$wsTasks.Cells.Item(11, 11).Value = 0
$wsTasks.Cells.Item(11, 12).Value = " [SECRET, EQUAL, 75lkfq7psmamjg9q65xdy]"
$wsTasks.Cells.Item(11, 13).Value = 1
$wsTasks.Cells.Item(11, 14).Value = "GameBus Studio"
